$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - (Intercept)
$ws.Range("B2").Value = 22869.23127
$ws.Range("D2").Value = 474.965166

# Row 3 - household_group_collapsed
$ws.Range("B3").Value = 1325.836395
$ws.Range("D3").Value = 13.767977
$ws.Range("E3").Value = 0.000002

# Row 4 - Residuals
$ws.Range("B4").Value = 10737.28968
$ws.Range("C4").Value = 223

# Row 5 - SM-Control
$ws.Range("G5").Value = -2.64401
$ws.Range("H5").Value = -5.674508
$ws.Range("I5").Value = 0.386488
$ws.Range("J5").Value = 0.100959

# Row 6 - SM + Traps-Control
$ws.Range("G6").Value = 2.820979
$ws.Range("H6").Value = -0.442226
$ws.Range("I6").Value = 6.084184
$ws.Range("J6").Value = 0.10519

# Row 7 - SM + Traps-SM
$ws.Range("G7").Value = 5.464989
$ws.Range("H7").Value = 2.997972
$ws.Range("I7").Value = 7.932006
$ws.Range("J7").Value = 0.000001
